# Update ServiceRequest rows to use a single required skill instead of a
# comma-separated list of multiple skills (required_skills M2M -> required_skill FK).
# Also update the ServiceType text to match the single skill, and trim the
# Technician "Skills" demo data down to a single skill per the sample data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Customer john_doe: RequiredSkills (J) and ServiceType (L)
$ws.Range("J2").Value = "Electric"
$ws.Range("L2").Value = "Electrical Repair"

# Row 3 - Customer jane_wilson: RequiredSkills (J) and ServiceType (L)
$ws.Range("J3").Value = "HVAC"
$ws.Range("L3").Value = "HVAC Maintenance"

# Row 4 - Customer robert_smith: RequiredSkills (J)
$ws.Range("J4").Value = "Plumbing"

# Row 5 - Technician tech_marie: Skills (Q)
$ws.Range("Q5").Value = "Electric"

# Row 6 - Technician tech_paul: Skills (Q)
$ws.Range("Q6").Value = "HVAC"

# Row 7 - Technician tech_sarah: Skills (Q)
$ws.Range("Q7").Value = "Plumbing"
